$d = $word.ActiveDocument
$d.Content.Find.Execute("python-powered", $true, $false, $false, $false, $false,
                         $true, 1, $false, "black", 2)
